$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the row for "Warehouse Employee / be able to update the stock of the
# inventory. / I can keep the products up to date." (sprint-2 backlog item
# that was completed / removed from the backlog).
$ws.Rows(11).Delete()

# Re-prioritize the remaining backlog items (column A) for sprint 2.
$ws.Range("A2").Value = 0   # Member / register for an account
$ws.Range("A3").Value = 0   # Member / browse items
$ws.Range("A5").Value = 0   # Manager / add and remove employees
$ws.Range("A6").Value = 0   # Member / borrow items
$ws.Range("A7").Value = 1   # Member / view my history
$ws.Range("A8").Value = 1   # Warehouse employee / view detailed list of outgoing orders
$ws.Range("A9").Value = 1   # Manager / manage the stock
$ws.Range("A11").Value = 1  # Manager / view records of items tracked by which employee
$ws.Range("A12").Value = 1  # Manager / track when things are received and shipped

# Re-sort the backlog by Priority (column A), stable on ties, same as the
# author re-running the worksheet's sort after editing priorities.
$ws.Range("A2:D20").Sort($ws.Range("A2:A20"))

# Leave the selection where the author left it after the edit.
$ws.Range("A11").Select()
